$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Header: volume/number line and week-covering date line (rich text -> plain concatenated text)
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"


# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '#,##0'
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(14, 5).Value = -100
$ws.Cells.Item(14, 7).NumberFormat = '#,##0'
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(14, 8).Value = -100
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(14, 14).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(14, 14).Value = -100

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0'
$ws.Cells.Item(15, 5).Value = '***.*'
$ws.Cells.Item(15, 12).Value = 150
$ws.Cells.Item(15, 13).Value = 150

# Row 16
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(16, 5).Value = 100
$ws.Cells.Item(16, 6).Value = 11
$ws.Cells.Item(16, 8).Value = 175
$ws.Cells.Item(16, 9).Value = 23
$ws.Cells.Item(16, 10).Value = 22
$ws.Cells.Item(16, 11).Value = 4.545454545454
$ws.Cells.Item(16, 12).Value = -11.538461538461
$ws.Cells.Item(16, 13).Value = -43.90243902439
$ws.Cells.Item(16, 14).Value = -81.6

# Row 17
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 8
$ws.Cells.Item(17, 5).Value = -75
$ws.Cells.Item(17, 6).Value = 17
$ws.Cells.Item(17, 7).Value = 16
$ws.Cells.Item(17, 8).Value = 6.25
$ws.Cells.Item(17, 9).Value = 43
$ws.Cells.Item(17, 10).Value = 33
$ws.Cells.Item(17, 11).Value = 30.30303030303
$ws.Cells.Item(17, 12).Value = -18.867924528301
$ws.Cells.Item(17, 13).Value = 26.470588235294
$ws.Cells.Item(17, 14).Value = 10.25641025641

# Row 18
$ws.Cells.Item(18, 3).NumberFormat = '#,##0'
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).NumberFormat = '#,##0'
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(18, 5).Value = 100
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(18, 7).Value = 3
$ws.Cells.Item(18, 8).Value = 33.333333333333
$ws.Cells.Item(18, 9).Value = 6
$ws.Cells.Item(18, 10).Value = 10
$ws.Cells.Item(18, 11).Value = -40
$ws.Cells.Item(18, 12).Value = -53.846153846153
$ws.Cells.Item(18, 13).Value = -87.755102040816
$ws.Cells.Item(18, 14).Value = -96

# Row 19
$ws.Cells.Item(19, 3).Value = 5
$ws.Cells.Item(19, 4).Value = 5
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 22
$ws.Cells.Item(19, 7).Value = 26
$ws.Cells.Item(19, 8).Value = -15.384615384615
$ws.Cells.Item(19, 9).Value = 47
$ws.Cells.Item(19, 10).Value = 74
$ws.Cells.Item(19, 11).Value = -36.486486486486
$ws.Cells.Item(19, 12).Value = -9.615384615384
$ws.Cells.Item(19, 13).Value = -2.083333333333
$ws.Cells.Item(19, 14).Value = -44.705882352941

# Row 20
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 200
$ws.Cells.Item(20, 6).Value = 10
$ws.Cells.Item(20, 7).Value = 5
$ws.Cells.Item(20, 8).Value = 100
$ws.Cells.Item(20, 9).Value = 25
$ws.Cells.Item(20, 10).Value = 29
$ws.Cells.Item(20, 11).Value = -13.793103448275
$ws.Cells.Item(20, 12).Value = -26.470588235294
$ws.Cells.Item(20, 13).Value = 4.166666666666
$ws.Cells.Item(20, 14).Value = -94.959677419354

# Row 21
$ws.Cells.Item(21, 3).Value = 16
$ws.Cells.Item(21, 4).Value = 18
$ws.Cells.Item(21, 5).Value = -11.111111111111
$ws.Cells.Item(21, 6).Value = 64
$ws.Cells.Item(21, 7).Value = 57
$ws.Cells.Item(21, 8).Value = 12.280701754386
$ws.Cells.Item(21, 9).Value = 149
$ws.Cells.Item(21, 10).Value = 172
$ws.Cells.Item(21, 11).Value = -13.372093023255
$ws.Cells.Item(21, 12).Value = -18.131868131868
$ws.Cells.Item(21, 13).Value = -24.747474747474
$ws.Cells.Item(21, 14).Value = -83.499446290144

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '#,##0'
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(22, 5).Value = -100
$ws.Cells.Item(22, 6).NumberFormat = '@'
$ws.Cells.Item(22, 6).Value = '0'
$ws.Cells.Item(22, 7).NumberFormat = '#,##0'
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(22, 8).Value = -100
$ws.Cells.Item(22, 10).Value = 5
$ws.Cells.Item(22, 11).Value = -80

# Row 23
$ws.Cells.Item(23, 6).Value = 8
$ws.Cells.Item(23, 7).Value = 4
$ws.Cells.Item(23, 8).Value = 100
$ws.Cells.Item(23, 9).Value = 23
$ws.Cells.Item(23, 10).Value = 21
$ws.Cells.Item(23, 11).Value = 9.523809523809
$ws.Cells.Item(23, 12).Value = -8
$ws.Cells.Item(23, 13).Value = 228.571428571429

# Row 24
$ws.Cells.Item(24, 3).Value = 24
$ws.Cells.Item(24, 4).Value = 10
$ws.Cells.Item(24, 5).Value = 140
$ws.Cells.Item(24, 6).Value = 55
$ws.Cells.Item(24, 7).Value = 37
$ws.Cells.Item(24, 8).Value = 48.648648648648
$ws.Cells.Item(24, 9).Value = 121
$ws.Cells.Item(24, 10).Value = 114
$ws.Cells.Item(24, 11).Value = 6.140350877192
$ws.Cells.Item(24, 12).Value = -17.123287671232
$ws.Cells.Item(24, 13).Value = 30.10752688172

# Row 25
$ws.Cells.Item(25, 3).Value = 7
$ws.Cells.Item(25, 5).Value = 600
$ws.Cells.Item(25, 6).Value = 10
$ws.Cells.Item(25, 7).Value = 5
$ws.Cells.Item(25, 8).Value = 100
$ws.Cells.Item(25, 9).Value = 30
$ws.Cells.Item(25, 10).Value = 27
$ws.Cells.Item(25, 11).Value = 11.111111111111
$ws.Cells.Item(25, 12).Value = -9.090909090909

# Row 26
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 5).Value = -87.5
$ws.Cells.Item(26, 6).Value = 27
$ws.Cells.Item(26, 7).Value = 24
$ws.Cells.Item(26, 8).Value = 12.5
$ws.Cells.Item(26, 9).Value = 75
$ws.Cells.Item(26, 10).Value = 75
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = -3.846153846153
$ws.Cells.Item(26, 13).Value = -17.582417582417

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0'
$ws.Cells.Item(27, 5).Value = '***.*'
$ws.Cells.Item(27, 12).Value = 100

# Row 28
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 8).Value = -66.666666666666
$ws.Cells.Item(28, 9).Value = 10
$ws.Cells.Item(28, 10).Value = 9
$ws.Cells.Item(28, 11).Value = 11.111111111111
$ws.Cells.Item(28, 12).Value = 150

# Row 29
$ws.Cells.Item(29, 3).NumberFormat = '#,##0'
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).NumberFormat = '#,##0'
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).NumberFormat = '#,##0'
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = -50
$ws.Cells.Item(29, 9).Value = 2
$ws.Cells.Item(29, 10).Value = 4
$ws.Cells.Item(29, 11).Value = -50
$ws.Cells.Item(29, 12).Value = -50
$ws.Cells.Item(29, 13).Value = -33.333333333333
$ws.Cells.Item(29, 14).Value = -66.666666666666

# Row 30
$ws.Cells.Item(30, 3).NumberFormat = '#,##0'
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).NumberFormat = '#,##0'
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).NumberFormat = '#,##0'
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = 2
$ws.Cells.Item(30, 8).Value = -50
$ws.Cells.Item(30, 9).Value = 2
$ws.Cells.Item(30, 10).Value = 4
$ws.Cells.Item(30, 11).Value = -50
$ws.Cells.Item(30, 12).Value = -50
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(30, 14).Value = -66.666666666666
